$d = $word.ActiveDocument

# Add a new bullet-list item after the last paragraph ("Wiggle - maybe a
# better result?"), inheriting the same list/paragraph formatting.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "Title?"
